$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "65.158.92"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "3.531.16"
$ws.Range("E3").Value = "  +2.51%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue "D5" "599.44"
$ws.Range("E5").Value = "  +1.13%  "
Set-TextValue "D6" "138.68"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "3.530.82"
$ws.Range("E7").Value = "  +2.55%  "
$ws.Range("E8").Value = "  +0.02%  "
Set-TextValue "D9" "0.490"
$ws.Range("E9").Value = "  -3.31%  "
$ws.Range("E10").Value = "  +1.21%  "
Set-TextValue "D11" "6.91"
$ws.Range("E11").Value = "  -6.09%  "
Set-TextValue "D12" "0.389"
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("D13").Value = "4.130.78"
$ws.Range("E13").Value = "  +2.53%  "
Set-TextValue "D14" "0.0000184"
$ws.Range("E14").Value = "  +1.18%  "
Set-TextValue "D15" "27.08"
$ws.Range("D16").Value = "3.536.94"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "65.235.54"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("E19").Value = "  +4.40%  "
Set-TextValue "D20" "5.93"
$ws.Range("E20").Value = "  +0.09%  "
Set-TextValue "D21" "14.28"
$ws.Range("E21").Value = "  +3.46%  "
Set-TextValue "D22" "392.47"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("E23").Value = "  +2.85%  "
$ws.Range("D24").Value = "3.671.37"
$ws.Range("E24").Value = "  +2.27%  "
Set-TextValue "D25" "73.71"
$ws.Range("E25").Value = "  +0.07%  "
Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +6.82%  "
Set-TextValue "D28" "7.64"
$ws.Range("E28").Value = "  +5.76%  "
Set-TextValue "D29" "0.999"
$ws.Range("E29").Value = "  -0.10%  "
Set-TextValue "D30" "2.28"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("E31").Value = "  -1.78%  "
$ws.Range("D32").Value = "3.541.88"
$ws.Range("E32").Value = "  +2.55%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  +2.83%  "
$ws.Range("E35").Value = "  -0.37%  "
Set-TextValue "D36" "1.26"
$ws.Range("E36").Value = "  +5.63%  "
Set-TextValue "D37" "6.94"
$ws.Range("E37").Value = "  -0.49%  "
Set-TextValue "D38" "168.59"
$ws.Range("E38").Value = "  -2.56%  "
Set-TextValue "D39" "1.54"
$ws.Range("E39").Value = "  +3.80%  "
$ws.Range("E40").Value = "  +2.65%  "
Set-TextValue "D41" "0.0802"
$ws.Range("E41").Value = "  +4.22%  "
$ws.Range("E42").Value = "  -0.55%  "
Set-TextValue "D43" "26.21"
$ws.Range("E43").Value = "  +12.90%  "
Set-TextValue "D44" "42.75"
$ws.Range("E44").Value = "  -2.64%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  +2.00%  "
Set-TextValue "D48" "1.17"
$ws.Range("E48").Value = "  +4.46%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D49" "6.79"
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.405.87"
$ws.Range("E50").Value = "  +8.67%  "
Set-TextValue "D51" "302.06"
$ws.Range("E51").Value = "  +6.80%  "
